# Update "想去人数" (F column) and one "最低票价" (G43 on 展览) counts
# to reflect the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 70
$ws.Range("F5").Value  = 96
$ws.Range("F6").Value  = 499
$ws.Range("F7").Value  = 4859
$ws.Range("F8").Value  = 4859
$ws.Range("F13").Value = 205
$ws.Range("F15").Value = 7849
$ws.Range("F20").Value = 1942
$ws.Range("F21").Value = 1942
$ws.Range("F25").Value = 2101
$ws.Range("F27").Value = 3
$ws.Range("F28").Value = 6260
$ws.Range("F30").Value = 48
$ws.Range("F34").Value = 6630
$ws.Range("F35").Value = 1
$ws.Range("F41").Value = 22
$ws.Range("G43").Value = 55
$ws.Range("F44").Value = 2486
$ws.Range("F48").Value = 48
$ws.Range("F49").Value = 466

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 25

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 70
$ws.Range("F7").Value  = 96
$ws.Range("F10").Value = 499
$ws.Range("F11").Value = 4859
$ws.Range("F12").Value = 4859
$ws.Range("F17").Value = 7849
$ws.Range("F21").Value = 1942
$ws.Range("F27").Value = 2101
$ws.Range("F31").Value = 3
$ws.Range("F32").Value = 6260
$ws.Range("F35").Value = 48
$ws.Range("F37").Value = 6630
$ws.Range("F41").Value = 22
$ws.Range("F47").Value = 48
$ws.Range("F48").Value = 466
$ws.Range("F51").Value = 25
